# Edit: 
#  1) Change the table style on the slide 5 table from the custom
#     "Table_0" style to the built-in Medium Style 2 - Accent 1 style.
#  2) Re-colour the presentation's theme (ppt/theme/theme1.xml, the theme
#     that slideMaster1.xml / the slides actually use) from the old
#     "Integral" / "Red Violet" palette to the standard Office theme
#     palette - i.e. applying the "Office Theme" design to the deck.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{E0FC41E1-E40E-453F-AC60-5DA536CC9BDA}")
        }
    }
}

# --- 2) Theme colours -------------------------------------------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$theme  = $master.Theme
$colors = $theme.ThemeColorScheme

# Office theme colour scheme, in clrScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeRGB = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $hexColor = $officeRGB[$i - 1]
    $r = [math]::Floor($hexColor / 0x10000) % 0x100
    $g = [math]::Floor($hexColor / 0x100) % 0x100
    $b = $hexColor % 0x100
    $bgr = $b * 0x10000 + $g * 0x100 + $r
    $colors.Item($i).RGB = $bgr
}
